# Update the "Data Type on Computer" column to "Data Type in R",
# and replace the old "String"/"Logical" values with "Character".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# Header
$ws.Range("C1").Value = "Data Type in R"

# Body values: String / Logical -> Character (Double stays Double)
$ws.Range("C3").Value = "Character"
$ws.Range("C4").Value = "Character"
$ws.Range("C7").Value = "Character"
$ws.Range("C11").Value = "Character"
$ws.Range("C12").Value = "Character"
